# Refresh Adam15-Itgav ligand/receptor expression stats with the latest TPM run.
# Sending-cluster ligand columns (G:J), target-cluster receptor columns (M:P) and the
# derived edge-weight columns (Q:T) are recomputed per row from the updated TPM table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 62.74008966666666
$ws.Range("H2").Value = 188.220269
$ws.Range("I2").Value = 0.6543216524118473
$ws.Range("J2").Value = 0.6543216524118471
$ws.Range("M2").Value = 8.820647333333334
$ws.Range("N2").Value = 26.461942
$ws.Range("O2").Value = 0.06415146660411865
$ws.Range("P2").Value = 0.06415146660411865
$ws.Range("Q2").Value = 553.4082046113774
$ws.Range("R2").Value = 4980.673841502397
$ws.Range("S2").Value = 0.04197569363305035
$ws.Range("T2").Value = 0.04197569363305034

# Row 3
$ws.Range("G3").Value = 62.74008966666666
$ws.Range("H3").Value = 188.220269
$ws.Range("I3").Value = 0.6543216524118473
$ws.Range("J3").Value = 0.6543216524118471
$ws.Range("O3").Value = 0.3979101621202897
$ws.Range("P3").Value = 0.3979101621202898
$ws.Range("Q3").Value = 3432.606611700971
$ws.Range("R3").Value = 30893.45950530874
$ws.Range("S3").Value = 0.260361234790014
$ws.Range("T3").Value = 0.260361234790014

# Row 4
$ws.Range("G4").Value = 62.74008966666666
$ws.Range("H4").Value = 188.220269
$ws.Range("I4").Value = 0.6543216524118473
$ws.Range("J4").Value = 0.6543216524118471
$ws.Range("M4").Value = 21.90816333333333
$ws.Range("N4").Value = 65.72449
$ws.Range("O4").Value = 0.1593353362087987
$ws.Range("P4").Value = 0.1593353362087987
$ws.Range("Q4").Value = 1374.520131965312
$ws.Range("R4").Value = 12370.68118768781
$ws.Range("S4").Value = 0.1042565604757384
$ws.Range("T4").Value = 0.1042565604757384

# Row 5
$ws.Range("G5").Value = 62.74008966666666
$ws.Range("H5").Value = 188.220269
$ws.Range("I5").Value = 0.6543216524118473
$ws.Range("J5").Value = 0.6543216524118471
$ws.Range("M5").Value = 52.056859
$ws.Range("N5").Value = 156.170577
$ws.Range("O5").Value = 0.3786030350667928
$ws.Range("P5").Value = 0.3786030350667929
$ws.Range("Q5").Value = 3266.052001425023
$ws.Range("R5").Value = 29394.4680128252
$ws.Range("S5").Value = 0.2477281635130444
$ws.Range("T5").Value = 0.2477281635130444

# Row 6
$ws.Range("I6").Value = 0.1782000513806195
$ws.Range("J6").Value = 0.1782000513806195
$ws.Range("M6").Value = 8.820647333333334
$ws.Range("N6").Value = 26.461942
$ws.Range("O6").Value = 0.06415146660411865
$ws.Range("P6").Value = 0.06415146660411865
$ws.Range("Q6").Value = 150.7169602789355
$ws.Range("R6").Value = 1356.45264251042
$ws.Range("S6").Value = 0.01143179464499604
$ws.Range("T6").Value = 0.01143179464499604

# Row 7
$ws.Range("I7").Value = 0.1782000513806195
$ws.Range("J7").Value = 0.1782000513806195
$ws.Range("O7").Value = 0.3979101621202897
$ws.Range("P7").Value = 0.3979101621202898
$ws.Range("S7").Value = 0.07090761133470629
$ws.Range("T7").Value = 0.07090761133470627

# Row 8
$ws.Range("I8").Value = 0.1782000513806195
$ws.Range("J8").Value = 0.1782000513806195
$ws.Range("M8").Value = 21.90816333333333
$ws.Range("N8").Value = 65.72449
$ws.Range("O8").Value = 0.1593353362087987
$ws.Range("P8").Value = 0.1593353362087987
$ws.Range("Q8").Value = 374.3412085433222
$ws.Range("R8").Value = 3369.0708768899
$ws.Range("S8").Value = 0.02839356509915621
$ws.Range("T8").Value = 0.02839356509915621

# Row 9
$ws.Range("I9").Value = 0.1782000513806195
$ws.Range("J9").Value = 0.1782000513806195
$ws.Range("M9").Value = 52.056859
$ws.Range("N9").Value = 156.170577
$ws.Range("O9").Value = 0.3786030350667928
$ws.Range("P9").Value = 0.3786030350667929
$ws.Range("Q9").Value = 889.4870471126966
$ws.Range("R9").Value = 8005.38342401427
$ws.Range("S9").Value = 0.06746708030176098
$ws.Range("T9").Value = 0.06746708030176098

# Row 10
$ws.Range("G10").Value = 2.950144666666667
$ws.Range("H10").Value = 8.850434
$ws.Range("I10").Value = 0.03076730593473967
$ws.Range("J10").Value = 0.03076730593473966
$ws.Range("M10").Value = 8.820647333333334
$ws.Range("N10").Value = 26.461942
$ws.Range("O10").Value = 0.06415146660411865
$ws.Range("P10").Value = 0.06415146660411865
$ws.Range("Q10").Value = 26.02218568698089
$ws.Range("R10").Value = 234.199671182828
$ws.Range("S10").Value = 0.001973767799171153
$ws.Range("T10").Value = 0.001973767799171153

# Row 11
$ws.Range("G11").Value = 2.950144666666667
$ws.Range("H11").Value = 8.850434
$ws.Range("I11").Value = 0.03076730593473967
$ws.Range("J11").Value = 0.03076730593473966
$ws.Range("O11").Value = 0.3979101621202897
$ws.Range("P11").Value = 0.3979101621202898
$ws.Range("Q11").Value = 161.4069431853967
$ws.Range("R11").Value = 1452.66248866857
$ws.Range("S11").Value = 0.01224262369249681
$ws.Range("T11").Value = 0.01224262369249681

# Row 12
$ws.Range("G12").Value = 2.950144666666667
$ws.Range("H12").Value = 8.850434
$ws.Range("I12").Value = 0.03076730593473967
$ws.Range("J12").Value = 0.03076730593473966
$ws.Range("M12").Value = 21.90816333333333
$ws.Range("N12").Value = 65.72449
$ws.Range("O12").Value = 0.1593353362087987
$ws.Range("P12").Value = 0.1593353362087987
$ws.Range("Q12").Value = 64.63225121429556
$ws.Range("R12").Value = 581.69026092866
$ws.Range("S12").Value = 0.004902319035350712
$ws.Range("T12").Value = 0.00490231903535071

# Row 13
$ws.Range("G13").Value = 2.950144666666667
$ws.Range("H13").Value = 8.850434
$ws.Range("I13").Value = 0.03076730593473967
$ws.Range("J13").Value = 0.03076730593473966
$ws.Range("M13").Value = 52.056859
$ws.Range("N13").Value = 156.170577
$ws.Range("O13").Value = 0.3786030350667928
$ws.Range("P13").Value = 0.3786030350667929
$ws.Range("Q13").Value = 153.5752649422687
$ws.Range("R13").Value = 1382.177384480418
$ws.Range("S13").Value = 0.01164859540772098
$ws.Range("T13").Value = 0.01164859540772098

# Row 14
$ws.Range("G14").Value = 13.108629
$ws.Range("H14").Value = 39.325887
$ws.Range("I14").Value = 0.1367109902727936
$ws.Range("J14").Value = 0.1367109902727935
$ws.Range("M14").Value = 8.820647333333334
$ws.Range("N14").Value = 26.461942
$ws.Range("O14").Value = 0.06415146660411865
$ws.Range("P14").Value = 0.06415146660411865
$ws.Range("Q14").Value = 115.626593432506
$ws.Range("R14").Value = 1040.639340892554
$ws.Range("S14").Value = 0.008770210526901107
$ws.Range("T14").Value = 0.008770210526901105

# Row 15
$ws.Range("G15").Value = 13.108629
$ws.Range("H15").Value = 39.325887
$ws.Range("I15").Value = 0.1367109902727936
$ws.Range("J15").Value = 0.1367109902727935
$ws.Range("O15").Value = 0.3979101621202897
$ws.Range("P15").Value = 0.3979101621202898
$ws.Range("Q15").Value = 717.193214335515
$ws.Range("R15").Value = 6454.738929019635
$ws.Range("S15").Value = 0.05439869230307264
$ws.Range("T15").Value = 0.05439869230307264

# Row 16
$ws.Range("G16").Value = 13.108629
$ws.Range("H16").Value = 39.325887
$ws.Range("I16").Value = 0.1367109902727936
$ws.Range("J16").Value = 0.1367109902727935
$ws.Range("M16").Value = 21.90816333333333
$ws.Range("N16").Value = 65.72449
$ws.Range("O16").Value = 0.1593353362087987
$ws.Range("P16").Value = 0.1593353362087987
$ws.Range("Q16").Value = 287.18598520807
$ws.Range("R16").Value = 2584.67386687263
$ws.Range("S16").Value = 0.02178289159855337
$ws.Range("T16").Value = 0.02178289159855336

# Row 17
$ws.Range("G17").Value = 13.108629
$ws.Range("H17").Value = 39.325887
$ws.Range("I17").Value = 0.1367109902727936
$ws.Range("J17").Value = 0.1367109902727935
$ws.Range("M17").Value = 52.056859
$ws.Range("N17").Value = 156.170577
$ws.Range("O17").Value = 0.3786030350667928
$ws.Range("P17").Value = 0.3786030350667929
$ws.Range("Q17").Value = 682.3940515363109
$ws.Range("R17").Value = 6141.546463826799
$ws.Range("S17").Value = 0.05175919584426644
$ws.Range("T17").Value = 0.05175919584426644
